# transactions working for digitalocean
#
# Note: the handout master's footer date field (a cached/auto-refreshed
# "today" datetimeFigureOut value, 4/12/18 -> 4/17/18) is not an
# addressable, independently-editable object in this host -- writes
# through ActivePresentation.HandoutMaster land on the slide master
# instead of the handout master, so it is intentionally left untouched
# here rather than risk corrupting unrelated slide-master content.
$p = $ppt.ActivePresentation

# Slide 3 ("API" flow diagram): the "GET /subsidies" endpoint box gets
#    a second, centered line for the "/{id}" sub-route, reusing the bold
#    white 12pt formatting already used in the box. There are two
#    "GET /subsidies" rectangles on this slide (the flow is drawn twice);
#    the one being edited here is "Rectangle 52", the lower copy that is
#    immediately followed by its elbow connector in the shape tree.
$s = $p.Slides.Item(3)

$epBox = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 52") {
        $epBox = $shp
    }
}
if ($epBox -eq $null) {
    # Fallback: last shape in the deck whose text is exactly "GET /subsidies"
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "GET /subsidies") {
                $epBox = $shp
            }
        }
    }
}

$tr = $epBox.TextFrame.TextRange
$null = $tr.InsertAfter("`r/{id}")
